$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 112, shifting existing rows 112:214 down to 113:215.
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with this week's data (the rest of the
# columns carry the same values as the row that used to occupy position 112,
# now at 113 - only the columns below actually differ).
$ws.Cells.Item(112, 1).Value2 = 4
$ws.Cells.Item(112, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(112, 3).Value2 = "Los Lagos"
$ws.Cells.Item(112, 4).Value2 = 44586
$ws.Cells.Item(112, 5).Value2 = 10
$ws.Cells.Item(112, 6).Value2 = 100112043
$ws.Cells.Item(112, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(112, 8).Value2 = "Sin especificar"
$ws.Cells.Item(112, 9).Value2 = "Primera"
$ws.Cells.Item(112, 10).Value2 = 400
$ws.Cells.Item(112, 11).Value2 = 12000
$ws.Cells.Item(112, 12).Value2 = 12500
$ws.Cells.Item(112, 13).Value2 = 12250
$ws.Cells.Item(112, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(112, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(112, 16).Value2 = 204
$ws.Cells.Item(112, 17).Value2 = 60
$ws.Cells.Item(112, 18).Value2 = "Hortaliza"
